$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the support document label into D6
$ws.Range("D6").Value = "MFJ-1278"

# Reflect the last active cell/selection used when the file was saved
$ws.Range("C12").Select()
